# feat: add 2022-Q1 data
#
# 1. Duplicate the "2021-Q4" sheet (it already carries the right header/style
#    layout for a per-fund holdings sheet) and place the copy right before
#    "总计", then rename it "2022-Q1" and update its figures.
# 2. Insert the new "2022-Q1" row at the top of the "总计" sheet's data,
#    pushing the existing rows down.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Step 1: new "2022-Q1" sheet
# ---------------------------------------------------------------------
$q4 = $wb.Worksheets.Item("2021-Q4")
$total = $wb.Worksheets.Item("总计")

$q4.Copy($total)

# Re-fetch "总计" - its Index moved now that a sheet was inserted before it.
$total = $wb.Worksheets.Item("总计")
$q1_2022 = $wb.Worksheets.Item($total.Index - 1)
$q1_2022.Name = "2022-Q1"

# Row 2 - fund 010343
$q1_2022.Range("D2").Value = "'0.20"
$q1_2022.Range("E2").Value = "'93.65"
$q1_2022.Range("F2").Value = "'4.19"
$q1_2022.Range("G2").Value = "'0.0084"
$q1_2022.Range("H2").Value = 5

# Row 3 - fund 010344 (D3 "0.06" is unchanged from the copied "2021-Q4" sheet)
$q1_2022.Range("E3").Value = "'93.65"
$q1_2022.Range("F3").Value = "'4.19"
$q1_2022.Range("G3").Value = "'0.0025"
$q1_2022.Range("H3").Value = 5

# ---------------------------------------------------------------------
# Step 2: prepend a "2022-Q1" row into "总计"
# ---------------------------------------------------------------------

# Push the existing two rows down one slot first (row 3 <- row 2 data,
# row 4 <- row 3 data), then write the brand-new top row.
$total.Range("A4").Value = 2
$total.Range("B4").Value = "2021-Q1"
$total.Range("C4").Value = 2
$total.Range("D4").Value = 0.01

# Row 4's "A" cell is brand new (rows only went to 3 before) - pick up the
# index-column style (bold/centered/bordered) from row 3 above it.
$total.Range("A3").Copy()
$total.Range("A4").PasteSpecial(-4122)

$total.Range("A3").Value = 1
$total.Range("B3").Value = "2021-Q4"
$total.Range("C3").Value = 2
$total.Range("D3").Value = 0.02

$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q1"
$total.Range("C2").Value = 2
$total.Range("D2").Value = 0.01

# Restore the original active sheet/tab selection (the Copy() above made the
# new sheet active as a side effect).
$wb.Worksheets.Item("2021-Q1").Activate()
